# Hortaliza, Vega Monumental Concepción - Betarraga: insert a new weekly
# price pair (Primera/Segunda) at the top of the date-ordered data block
# (rows 404-529), pushing every subsequent pair down by 2 rows so the
# sheet grows from A1:R529 to A1:R531.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 404
$lastRow  = 529
$blockLen = $lastRow - $firstRow + 1      # 126 existing data rows to shift

# 1) Snapshot the existing D:P block (columns D..P) before touching anything.
#    Range.Value2 comes back as a 1-based [1..rows,1..cols] array.
$srcRange = $ws.Range("D$firstRow`:P$lastRow")
$srcVals  = $srcRange.Value2

# 2) Build the replacement block: 2 new rows on top, then the old block.
#    New-Object 'object[,]' produces a 0-based .NET array, and Range.Value2
#    assignment maps [0,0] to the top-left cell of the destination range -
#    so build this one 0-based.
$newBlockRows = $blockLen + 2             # 128 rows -> will land on 404..531
$newVals = New-Object 'object[,]' $newBlockRows,13

# columns within the D:P range, 0-based offsets:
# 0=D 1=E 2=F 3=G 4=H 5=I 6=J 7=K 8=L 9=M 10=N 11=O 12=P

# New row 0 (-> sheet row 404): Primera, fecha 45120
$newVals[0,0]  = 45120
$newVals[0,1]  = 8
$newVals[0,2]  = 100114014
$newVals[0,3]  = "Betarraga"
$newVals[0,4]  = "Sin especificar"
$newVals[0,5]  = "Primera"
$newVals[0,6]  = 200
$newVals[0,7]  = 700
$newVals[0,8]  = 700
$newVals[0,9]  = 700
$newVals[0,10] = "`$/paquete 5 unidades"
$newVals[0,11] = "Región Metropolitana"
$newVals[0,12] = 140

# New row 1 (-> sheet row 405): Segunda, fecha 45120
$newVals[1,0]  = 45120
$newVals[1,1]  = 8
$newVals[1,2]  = 100114014
$newVals[1,3]  = "Betarraga"
$newVals[1,4]  = "Sin especificar"
$newVals[1,5]  = "Segunda"
$newVals[1,6]  = 200
$newVals[1,7]  = 500
$newVals[1,8]  = 500
$newVals[1,9]  = 500
$newVals[1,10] = "`$/paquete 5 unidades"
$newVals[1,11] = "Región Metropolitana"
$newVals[1,12] = 100

# Remaining rows: the previously-existing block, shifted down by 2.
for ($r = 1; $r -le $blockLen; $r++) {
    for ($c = 1; $c -le 13; $c++) {
        $newVals[$r + 1, $c - 1] = $srcVals[$r, $c]
    }
}

# 3) Write the whole new block back out in one shot (rows 404..531).
$newLastRow = $firstRow + $newBlockRows - 1   # 531
$destRange = $ws.Range("D$firstRow`:P$newLastRow")
$destRange.Value2 = $newVals

# 4) The two brand-new sheet rows (530, 531) also need the columns that
#    sit outside D:P (A,B,C,Q,R) filled in - those are constant for every
#    data row in this sheet, so copy them from the row just above the
#    appended tail.
foreach ($col in @("A","B","C","Q","R")) {
    $ws.Range("$col$($lastRow)").Copy() | Out-Null
    $ws.Range("$col$($lastRow + 1)`:$col$newLastRow").PasteSpecial(-4104) | Out-Null
}
$excel.CutCopyMode = 0

# 5) Match the date-column number format (style index 2 in the original
#    file) on the two freshly-created D cells.
$ws.Range("D$($lastRow + 1)`:D$newLastRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"

"done: rows $firstRow..$newLastRow rewritten"
